$d = $word.ActiveDocument
$Wns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Helper: find the unique paragraph index whose trimmed text (stripping the
# trailing paragraph/cell marks) equals $text, optionally restricted to a
# given style name.
# ---------------------------------------------------------------------------
function Get-ParaIndexByText($text, $styleName) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text -replace "[\r\x07]", ""
        if ($t -eq $text) {
            if ($styleName -eq $null -or $p.Style.NameLocal -eq $styleName) {
                return $i
            }
        }
    }
    return -1
}

# ===========================================================================
# 1. "Time Series Analysis" paragraph -> bold
# ===========================================================================
$i = Get-ParaIndexByText "Time Series Analysis" $null
$p = $d.Paragraphs($i)
$xml = '<w:p xmlns:w="' + $Wns + '"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Time Series Analysis</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ===========================================================================
# 2. "Weather Modeling" .. "Air Inversion Forecasting" block: bold the three
#    headers, add "Search: " lines, the particulate-matter literature list
#    (5 hyperlinks + summary line), and a trailing "Search: " line.
# ===========================================================================
$f1 = $d.Content
$f1.Find.Execute("Weather Modeling", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $f1.Paragraphs(1).Range.Start

$f2 = $d.Content
$f2.Find.Execute("Air Inversion Forecasting", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $f2.Paragraphs(1).Range.End

$fullRange = $d.Range($startPos, $endPos)

$placeholder = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Weather Modeling</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Search: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Particulate Matter Forecasting</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Search: </w:t></w:r><w:r><w:t>particulate matter forecasting</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="90001" w:name="PMF_LINK_1"/><w:bookmarkEnd w:id="90001"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="90002" w:name="PMF_LINK_2"/><w:bookmarkEnd w:id="90002"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="90003" w:name="PMF_LINK_3"/><w:bookmarkEnd w:id="90003"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="90004" w:name="PMF_LINK_4"/><w:bookmarkEnd w:id="90004"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="90005" w:name="PMF_LINK_5"/><w:bookmarkEnd w:id="90005"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">There are dozens of papers on this topic which bodes well for an ML approach. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Air Inversion Forecasting</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Search: </w:t></w:r></w:p>
'@

$fullRange.InsertXML($placeholder)

# Fill in the five literature-review hyperlinks at the bookmark placeholders,
# each followed by a plain trailing space run, then remove the bookmark.
$links = @(
    @{ Name = "PMF_LINK_1"; Url = "https://www.sciencedirect.com/science/article/pii/S1877050920312060" },
    @{ Name = "PMF_LINK_2"; Url = "https://www.arl.noaa.gov/research/surface-atmosphere-exchange-home/o3-and-pm-2/" },
    @{ Name = "PMF_LINK_3"; Url = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC9723408/" },
    @{ Name = "PMF_LINK_4"; Url = "https://www.mdpi.com/2073-4433/13/9/1451" },
    @{ Name = "PMF_LINK_5"; Url = "https://ieeexplore.ieee.org/document/9359734" }
)

foreach ($link in $links) {
    $bm = $d.Bookmarks.Item($link.Name)
    $rng = $bm.Range
    $h = $d.Hyperlinks.Add($rng, $link.Url, "", "", $link.Url)
    $hEnd = $h.Range.End
    $spaceRng = $d.Range($hEnd, $hEnd)
    $spaceRng.InsertAfter(" ")
    $bm2 = $d.Bookmarks.Item($link.Name)
    $bm2.Delete()
}

# ===========================================================================
# 3. "Data Collection and Sourcing" heading -> add <w:lastRenderedPageBreak/>
# ===========================================================================
$i = Get-ParaIndexByText "Data Collection and Sourcing" "Heading 2"
$p = $d.Paragraphs($i)
$xml = '<w:p xmlns:w="' + $Wns + '"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="90010" w:name="_Data_Collection_and"/><w:bookmarkEnd w:id="90010"/><w:r><w:lastRenderedPageBreak/><w:t>Data Collection and Sourcing</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ===========================================================================
# 4. "F8" table cell -> remove <w:lastRenderedPageBreak/>
# ===========================================================================
$i = Get-ParaIndexByText "F8" $null
$p = $d.Paragraphs($i)
$xml = '<w:p xmlns:w="' + $Wns + '"><w:r><w:t>F8</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ===========================================================================
# 5. Wrap field-name cells with proofErr spellStart/spellEnd
# ===========================================================================
$fieldNames = @("windDIR", "windMPH", "precip", "mslp", "wxcodes")
foreach ($fn in $fieldNames) {
    $i = Get-ParaIndexByText $fn $null
    $p = $d.Paragraphs($i)
    $xml = '<w:p xmlns:w="' + $Wns + '"><w:proofErr w:type="spellStart"/><w:r><w:t>' + $fn + '</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
    $p.Range.InsertXML($xml)
}

# ===========================================================================
# 6. "F26" table cell -> add <w:lastRenderedPageBreak/>
# ===========================================================================
$i = Get-ParaIndexByText "F26" $null
$p = $d.Paragraphs($i)
$xml = '<w:p xmlns:w="' + $Wns + '"><w:r><w:lastRenderedPageBreak/><w:t>F26</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ===========================================================================
# 7. "Data Cleaning" (Heading 3) -> remove <w:lastRenderedPageBreak/>
# ===========================================================================
$i = Get-ParaIndexByText "Data Cleaning" "Heading 3"
$p = $d.Paragraphs($i)
$xml = '<w:p xmlns:w="' + $Wns + '"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Data Cleaning</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ===========================================================================
# 8. "Model Building and Evaluation" heading -> add <w:lastRenderedPageBreak/>
# ===========================================================================
$i = Get-ParaIndexByText "Model Building and Evaluation" "Heading 2"
$p = $d.Paragraphs($i)
$xml = '<w:p xmlns:w="' + $Wns + '"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="90011" w:name="_Model_Building_and"/><w:bookmarkEnd w:id="90011"/><w:r><w:lastRenderedPageBreak/><w:t>Model Building and Evaluation</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

Write-Output "All edits applied."
